$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 12 & 13: Cardano / Toncoin swap (with updated price + volume) ---
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.55"
$ws.Range("E12").Value = "  -5.68%  "

$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.310"
$ws.Range("E13").Value = "  -3.88%  "

# --- Rows 22 & 23: Dai / Uniswap swap (with updated price + volume) ---
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("E22").Value = "  -3.73%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.10%  "

# --- Remaining rows: price (D) / volume-1h (E) refresh ---
$ws.Range("D2").Value = "53.173.65"
$ws.Range("E2").Value = "  -5.32%  "
$ws.Range("D3").Value = "2.180.40"
$ws.Range("E3").Value = "  -7.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "479.72"
$ws.Range("E5").Value = "  -4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.14"
$ws.Range("E6").Value = "  -4.62%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  -5.57%  "
$ws.Range("D9").Value = "2.194.63"
$ws.Range("E9").Value = "  -7.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0907"
$ws.Range("E10").Value = "  -7.67%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D14").Value = "2.568.97"
$ws.Range("E14").Value = "  -7.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.81"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").Value = "53.080.56"
$ws.Range("E16").Value = "  -5.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000126"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").Value = "2.188.06"
$ws.Range("E18").Value = "  -7.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.46"
$ws.Range("E19").Value = "  -5.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.90"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "290.88"
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.08"
$ws.Range("E24").Value = "  -5.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.361"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").Value = "2.282.32"
$ws.Range("E27").Value = "  -7.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.143"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.93"
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.38"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.56"
$ws.Range("E32").Value = "  -4.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.995"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "0.0₃0652"
$ws.Range("E34").Value = "  -8.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.60"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.13"
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.13"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.805"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.62"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  -6.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.362"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.34"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "122.42"
$ws.Range("E45").Value = "  -5.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.67"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0870"
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.527"
$ws.Range("E48").Value = "  -6.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0464"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "226.09"
$ws.Range("E50").Value = "  -5.44%  "
$ws.Range("E51").Value = "  -4.46%  "
